$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '61.183.83'
$ws.Cells.Item(2, 5).Value = '  -1.51%  '
$ws.Cells.Item(3, 4).Value = '2.424.50'
$ws.Cells.Item(3, 5).Value = '  -0.85%  '
$ws.Cells.Item(4, 4).Value = '0.997'
$ws.Cells.Item(4, 5).Value = '  -0.10%  '
$ws.Cells.Item(5, 4).Value = '570.19'
$ws.Cells.Item(5, 5).Value = '  -2.03%  '
$ws.Cells.Item(6, 4).Value = '139.96'
$ws.Cells.Item(6, 5).Value = '  -1.92%  '
$ws.Cells.Item(7, 5).Value = '  +0.17%  '
$ws.Cells.Item(8, 4).Value = '0.528'
$ws.Cells.Item(8, 5).Value = '  -0.95%  '
$ws.Cells.Item(9, 4).Value = '2.410.08'
$ws.Cells.Item(9, 5).Value = '  -1.23%  '
$ws.Cells.Item(10, 5).Value = '  -2.03%  '
$ws.Cells.Item(11, 5).Value = '  -0.18%  '
$ws.Cells.Item(12, 5).Value = '  -2.75%  '
$ws.Cells.Item(14, 4).Value = '26.19'
$ws.Cells.Item(14, 5).Value = '  -0.69%  '
$ws.Cells.Item(15, 5).Value = '  -1.94%  '
$ws.Cells.Item(16, 4).Value = '2.840.96'
$ws.Cells.Item(16, 5).Value = '  -1.38%  '
$ws.Cells.Item(17, 4).Value = '61.054.26'
$ws.Cells.Item(17, 5).Value = '  -1.63%  '
$ws.Cells.Item(18, 4).Value = '2.402.25'
$ws.Cells.Item(18, 5).Value = '  -1.38%  '
$ws.Cells.Item(19, 4).Value = '7.94'
$ws.Cells.Item(19, 5).Value = '  +10.10%  '
$ws.Cells.Item(20, 5).Value = '  -0.53%  '
$ws.Cells.Item(21, 4).Value = '323.68'
$ws.Cells.Item(21, 5).Value = '  -0.76%  '
$ws.Cells.Item(22, 5).Value = '  -0.80%  '
$ws.Cells.Item(23, 5).Value = '  +2.16%  '
$ws.Cells.Item(24, 5).Value = '  +0.20%  '
$ws.Cells.Item(25, 5).Value = '  -3.85%  '
$ws.Cells.Item(26, 4).Value = '64.67'
$ws.Cells.Item(26, 5).Value = '  -1.39%  '
$ws.Cells.Item(27, 4).Value = '590.54'
$ws.Cells.Item(27, 5).Value = '  -1.27%  '
$ws.Cells.Item(28, 4).Value = '8.27'
$ws.Cells.Item(28, 5).Value = '  -9.36%  '
$ws.Cells.Item(30, 5).Value = '  -3.48%  '
$ws.Cells.Item(31, 4).Value = '7.90'
$ws.Cells.Item(31, 5).Value = '  -0.90%  '
$ws.Cells.Item(32, 4).Value = '1.35'
$ws.Cells.Item(32, 5).Value = '  -4.62%  '
$ws.Cells.Item(33, 5).Value = '  -3.75%  '
$ws.Cells.Item(34, 5).Value = '  -1.09%  '
$ws.Cells.Item(35, 5).Value = '  -0.09%  '
$ws.Cells.Item(36, 5).Value = '  -0.44%  '
$ws.Cells.Item(37, 4).Value = '4.62'
$ws.Cells.Item(37, 5).Value = '  -5.25%  '
$ws.Cells.Item(38, 4).Value = '151.83'
$ws.Cells.Item(38, 5).Value = '  -1.41%  '
$ws.Cells.Item(40, 4).Value = '18.24'
$ws.Cells.Item(41, 4).Value = '5.16'
$ws.Cells.Item(41, 5).Value = '  -2.41%  '
$ws.Cells.Item(42, 5).Value = '  +0.00%  '
$ws.Cells.Item(43, 5).Value = '  -1.92%  '
$ws.Cells.Item(44, 4).Value = '41.25'
$ws.Cells.Item(44, 5).Value = '  -4.68%  '
$ws.Cells.Item(45, 4).Value = '2.36'
$ws.Cells.Item(45, 5).Value = '  -6.94%  '
$ws.Cells.Item(46, 4).Value = '0.0₆0288'
$ws.Cells.Item(46, 5).Value = '  +4.35%  '
$ws.Cells.Item(47, 4).Value = '143.79'
$ws.Cells.Item(47, 5).Value = '  +1.44%  '
$ws.Cells.Item(48, 4).Value = '3.53'
$ws.Cells.Item(48, 5).Value = '  -2.39%  '
$ws.Cells.Item(49, 4).Value = '0.588'
$ws.Cells.Item(49, 5).Value = '  -2.19%  '
$ws.Cells.Item(50, 4).Value = '19.58'
$ws.Cells.Item(50, 5).Value = '  -1.18%  '
$ws.Cells.Item(51, 5).Value = '  -2.74%  '
